# Sender.py.docx edit: rewrite the UDP file-sender snippet to use a
# fixed-size header (file name + size) instead of separate struct-packed
# size + filename datagrams, add a socket timeout, rename several
# variables, and drop the old "end marker" datagram.

$d = $word.ActiveDocument

function Run-Xml([string]$text) {
    $esc = $text.Replace("&","&amp;").Replace("<","&lt;").Replace(">","&gt;")
    $out = "<w:r><w:t xml:space=`"preserve`">$esc</w:t></w:r>"
    return $out
}

function SpellRun-Xml([string]$text) {
    $esc = $text.Replace("&","&amp;").Replace("<","&lt;").Replace(">","&gt;")
    $out = "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>$esc</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>"
    return $out
}

function Set-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>$innerXml</w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
    $r.InsertXML($xml) | Out-Null
}

# Clears a freshly-inserted paragraph's placeholder empty <w:r></w:r> so it
# serializes as a truly empty paragraph (matching the source "<w:p/>" blank
# lines) instead of a paragraph holding one empty run.
function Clear-Para($paraIndex) {
    Set-ParaXml $paraIndex ""
}

# ---------------------------------------------------------------------
# Work from the bottom of the document upward so paragraph indices
# above the current edit point never shift under us.
# ---------------------------------------------------------------------

# 1. Remove the old "# Send end marker" / sock.sendto(b'__END__', ...) pair (paras 30-31).
$p1 = $d.Paragraphs(30)
$p2 = $d.Paragraphs(31)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$rng.Delete()

# 2. progress_callback(sent, file_size) -> progress_callback(sent_bytes, file_size)
$d.Paragraphs(28).Range.Text = "                progress_callback(sent_bytes, file_size)"

# 3. sent += len(data) -> sent_bytes += len(data), then a new blank line after it.
$a = Run-Xml "            sent_bytes += "
$b = SpellRun-Xml "len"
$c = Run-Xml "(data)"
$xml3 = $a + $b + $c
Set-ParaXml 26 $xml3
$d.Paragraphs(26).Range.InsertParagraphAfter() | Out-Null

# 4. data = f.read(CHUNK_SIZE) -> data = f.read(BUFFER_SIZE)
$d.Paragraphs(22).Range.Text = "            data = f.read(BUFFER_SIZE)"

# 5. with open(filename, 'rb') as f: -> with open(file_path, "rb") as f:
$a = Run-Xml "    with open(file_path, `""
$b = SpellRun-Xml "rb"
$c = Run-Xml "`") as f:"
$xml5 = $a + $b + $c
Set-ParaXml 20 $xml5

# 6. sent = 0 -> sent_bytes = 0
$d.Paragraphs(19).Range.Text = "    sent_bytes = 0"

# 7. # Send file content -> # Send file in chunks
$d.Paragraphs(18).Range.Text = "    # Send file in chunks"

# 8. Rework the header block (paras 14-16 -> 4 paragraphs).
$d.Paragraphs(14).Range.Text = "    # Send header"
$a = Run-Xml "    header = f`"{file_name}|{file_size}`".encode()."
$b = SpellRun-Xml "ljust"
$c = Run-Xml "(HEADER_SIZE, b'#')"
$xml8 = $a + $b + $c
Set-ParaXml 15 $xml8
$d.Paragraphs(16).Range.Text = "    sock.sendto(header, (ip, port))"
$d.Paragraphs(16).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(17).Range.Text = "    time.sleep(0.1)  # Give time for receiver to prepare"

# 9. "# Send filename" / base_name.../sock.sendto(base_name...) -> file_name=.. / file_size=..
$d.Paragraphs(10).Range.Delete()
$a = Run-Xml "    file_name = "
$b = SpellRun-Xml "os.path.basename"
$c = Run-Xml "(file_path)"
$xml9 = $a + $b + $c
Set-ParaXml 10 $xml9
$d.Paragraphs(11).Range.Delete()
$d.Paragraphs(10).Range.InsertParagraphAfter() | Out-Null
$a = Run-Xml "    file_size = "
$b = SpellRun-Xml "os.path.getsize"
$c = Run-Xml "(file_path)"
$xml9b = $a + $b + $c
Set-ParaXml 11 $xml9b

# 10. def send_file(filename, ...) -> def send_file(file_path, ...)
$d.Paragraphs(7).Range.Text = "def send_file(file_path, ip, port, progress_callback=None):"

# 11. New "sock.settimeout(2)" line right after the socket() call.
$d.Paragraphs(8).Range.InsertParagraphAfter() | Out-Null
$a = Run-Xml "    "
$b = SpellRun-Xml "sock.settimeout"
$c = Run-Xml "(2)"
$xml11 = $a + $b + $c
Set-ParaXml 9 $xml11

# 12. CHUNK_SIZE = 1024 -> BUFFER_SIZE = 4096, then new HEADER_SIZE = 1024 line.
$d.Paragraphs(5).Range.Text = "BUFFER_SIZE = 4096"
$d.Paragraphs(5).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(6).Range.Text = "HEADER_SIZE = 1024"

# 13. import struct -> import time
$d.Paragraphs(3).Range.Text = "import time"
